# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets
# F2: 2785 -> 2790
# F6: 1540 -> 1544
# F10: 87  -> 89
# F11: 14  -> 15

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 2790
    $ws.Range("F6").Value = 1544
    $ws.Range("F10").Value = 89
    $ws.Range("F11").Value = 15
}
